$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format for numeric-looking price values so Excel keeps them as text,
# matching the original inline-string (text) cell type.
$ws.Range("D2").Value = "29.756.19"

$ws.Range("D3").Value = "1.889.92"
$ws.Range("E3").Value = "  -0.99%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7934"
$ws.Range("E5").Value = "  -2.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.59"
$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3164"
$ws.Range("E8").Value = "  +1.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.58"
$ws.Range("E9").Value = "  -3.41%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07034"
$ws.Range("E10").Value = "  +0.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08047"
$ws.Range("E11").Value = "  +0.47%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7658"
$ws.Range("E12").Value = "  +2.83%  "

$ws.Range("D13").Value = "1.913.71"
$ws.Range("E13").Value = "  +0.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.315"
$ws.Range("E14").Value = "  +2.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.25"
$ws.Range("E15").Value = "  -0.32%  "

$ws.Range("D16").Value = "29.750.02"
$ws.Range("E16").Value = "  -0.73%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.79"
$ws.Range("E17").Value = "  -1.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.935"
$ws.Range("E18").Value = "  +1.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.78"
$ws.Range("E19").Value = "  -1.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007694"
$ws.Range("E20").Value = "  -1.21%  "

$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.08%  "

$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.165"
$ws.Range("E22").Value = "  +17.62%  "

$ws.Range("D23").Value = "2.149.45"
$ws.Range("E23").Value = "  -0.26%  "

$ws.Range("E24").Value = "  +0.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1623"
$ws.Range("E25").Value = "  +5.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.309"
$ws.Range("E26").Value = "  +1.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.92"
$ws.Range("E27").Value = "  -2.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.64"
$ws.Range("E28").Value = "  -1.18%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.056"
$ws.Range("E29").Value = "  -0.68%  "

$ws.Range("E30").Value = "  +1.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.534"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.438"
$ws.Range("E32").Value = "  +3.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05705"
$ws.Range("E33").Value = "  +3.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.077"
$ws.Range("E34").Value = "  +0.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.264"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7374"
$ws.Range("E36").Value = "  +0.97%  "

$ws.Range("E37").Value = "  +0.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.684"
$ws.Range("E38").Value = "  -0.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01911"
$ws.Range("E39").Value = "  -0.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.769"
$ws.Range("E40").Value = "  -0.68%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4401"
$ws.Range("E41").Value = "  -0.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.38"
$ws.Range("E42").Value = "  +0.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.838"
$ws.Range("E43").Value = "  -2.55%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8400"
$ws.Range("E44").Value = "  +0.25%  "

$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9997"
$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("D46").Value = "1.029.16"
$ws.Range("E46").Value = "  +4.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.16"
$ws.Range("E47").Value = "  +1.22%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.887"
$ws.Range("E48").Value = "  +1.75%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.850"
$ws.Range("E49").Value = "  -2.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.451"
$ws.Range("E50").Value = "  -1.62%  "

$ws.Range("D51").Value = "2.037.87"
$ws.Range("E51").Value = "  -1.02%  "
